# Fixing name of Sectors to be alligned with Baseline
#
# 1) Rename the 4 "Sector" header labels (row 3, columns D:G on every
#    yearly worksheet) from their short codes to their full Baseline names.
# 2) A handful of G7 ("Wires") totals carry a last-digit floating point
#    refresh that came along with the source data used to produce this
#    edit; apply those literal value corrections on the specific years
#    they touch.

$wb = $excel.ActiveWorkbook

$renames = @{
    "Nd" = "Neodymium"
    "Dy" = "Dysprosium"
    "Cu" = "Copper ores and concentrates"
    "Si" = "Raw silicon"
}

foreach ($ws in $wb.Worksheets) {
    if ($ws.Range("D3").Value2 -eq "Nd")  { $ws.Range("D3").Value2 = $renames["Nd"] }
    if ($ws.Range("E3").Value2 -eq "Dy")  { $ws.Range("E3").Value2 = $renames["Dy"] }
    if ($ws.Range("F3").Value2 -eq "Cu")  { $ws.Range("F3").Value2 = $renames["Cu"] }
    if ($ws.Range("G3").Value2 -eq "Si")  { $ws.Range("G3").Value2 = $renames["Si"] }
}

$g7Updates = @{
    "2010" = -18083.13819547712
    "2011" = -37637.75057429998
    "2013" = -143786.4232588847
    "2022" = -20743302.54629443
    "2023" = -37887279.27537362
    "2024" = -61530163.06617802
    "2026" = -90857625.99066542
    "2031" = -181665470.8205312
    "2035" = -255560285.057339
    "2040" = -412130944.2035097
    "2056" = -2861049668.1251
    "2067" = -395204002.9777587
    "2071" = -332106787.8604285
    "2072" = -312816625.1910507
    "2077" = -189667351.1309118
    "2080" = -213012274.4682261
    "2084" = -276735493.7468376
    "2086" = -293428808.0007818
}

foreach ($sheetName in $g7Updates.Keys) {
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $ws.Range("G7").Value2 = $g7Updates[$sheetName]
}
